# daily auto push: 2025-10-09 13:38 UTC
# Append the new day's log entry as the next row of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the first empty row below the existing data (column A has no gaps).
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# Column A holds the date as plain text (e.g. "2025/10/09"), matching the
# existing rows above it. Force a Text number format before assigning the
# value so Excel doesn't auto-convert the string into a date serial, then
# restore the default "Normal" style so the cell stays unformatted like its
# neighbours.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025/10/09"
$dateCell.Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = "木"
$ws.Cells.Item($newRow, 3).Value = 20
$ws.Cells.Item($newRow, 4).Value = 24
